$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows (34-36) following the existing pattern of data,
# continuing the usr_id / B column sequence (110033, 110034, 110035)
# with regcntr_id / machine_id = 10005.
$newRows = @(
    @(10005, 110033, 10005),
    @(10005, 110034, 10005),
    @(10005, 110035, 10005)
)

$r = 34
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
    $ws.Cells.Item($r, 8).Value = "now()"
    $r++
}

# Update the selection to match the post-edit state (cursor moved to the
# first empty row, with the rest of the column selected).
$ws.Range("A37:XFD1048576").Select()
